$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The organization's website changed (www.stat.kg -> www.stat.gov.kg) and the
# phone number lost its leading country code (+996 (0312) 32 46 36 -> (0312) 32 46 36).
# Update the website cell (B10) first, then the phone cell (B9), so that the
# shared-string table grows in the same order as the authored workbook.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B9").Value = "(0312) 32 46 36"

# Match the author's final on-screen selection (last cell they edited).
$null = $ws.Range("B9").Select()
